$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct mentor name typos ("Pref." -> "Prof." and remove the feminine ordinal sign)
$ws.Range("C4").Value = "Prof. Hirant Sanazar"
$ws.Range("C5").Value = "Prof. Ilza Nascimento Pintus"

# Add an (empty) underlined placeholder cell
$ws.Range("C6").Font.Underline = $true

# Update the active selection to C4
$ws.Range("C4").Select() | Out-Null
